$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.695.69"
$ws.Range("E2").Value = "  +2.14%  "

# Row 3
$ws.Range("D3").Value = "2.223.98"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'241.16"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'75.13"
$ws.Range("E7").Value = "  +3.24%  "

# Row 8
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("D9").Value = "'0.605"
$ws.Range("E9").Value = "  +2.76%  "

# Row 10
$ws.Range("D10").Value = "'41.52"
$ws.Range("E10").Value = "  +0.15%  "

# Row 11
$ws.Range("D11").Value = "'0.0932"
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("D12").Value = "'55.03"
$ws.Range("E12").Value = "  -2.21%  "

# Row 13
$ws.Range("D13").Value = "'6.92"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("E14").Value = "  -1.57%  "

# Row 15
$ws.Range("D15").Value = "2.555.88"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").Value = "'14.70"
$ws.Range("E16").Value = "  +3.98%  "

# Row 17
$ws.Range("D17").Value = "2.215.48"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("D18").Value = "'0.805"
$ws.Range("E18").Value = "  -2.57%  "

# Row 19
$ws.Range("D19").Value = "42.530.36"
$ws.Range("E19").Value = "  +2.00%  "

# Row 20
$ws.Range("E20").Value = "  +1.09%  "

# Row 21
$ws.Range("D21").Value = "'70.86"
$ws.Range("E21").Value = "  -1.10%  "

# Row 22
$ws.Range("D22").Value = "'5.96"
$ws.Range("E22").Value = "  -2.73%  "

# Row 23
$ws.Range("D23").Value = "'9.90"
$ws.Range("E23").Value = "  -7.96%  "

# Row 24
$ws.Range("D24").Value = "'229.89"
$ws.Range("E24").Value = "  +0.73%  "

# Row 25
$ws.Range("D25").Value = "'2.15"
$ws.Range("E25").Value = "  +6.15%  "

# Row 26
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").Value = "'10.96"
$ws.Range("E27").Value = "  -2.80%  "

# Row 28
$ws.Range("E28").Value = "  -7.08%  "

# Row 29
$ws.Range("E29").Value = "  -1.63%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.18"
$ws.Range("E30").Value = "  -0.82%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'172.86"
$ws.Range("E31").Value = "  +3.51%  "

# Row 32
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'36.31"
$ws.Range("E32").Value = "  +18.92%  "

# Row 33
$ws.Range("D33").Value = "'20.31"
$ws.Range("E33").Value = "  -0.17%  "

# Row 34
$ws.Range("D34").Value = "'0.0797"
$ws.Range("E34").Value = "  +0.97%  "

# Row 35
$ws.Range("E35").Value = "  -0.51%  "

# Row 36
$ws.Range("E36").Value = "  -0.94%  "

# Row 37
$ws.Range("E37").Value = "  -0.32%  "

# Row 38
$ws.Range("D38").Value = "'4.44"
$ws.Range("E38").Value = "  +4.88%  "

# Row 39
$ws.Range("E39").Value = "  +7.26%  "

# Row 40
$ws.Range("D40").Value = "'12.58"
$ws.Range("E40").Value = "  -2.58%  "

# Row 41
$ws.Range("E41").Value = "  +1.73%  "

# Row 42
$ws.Range("D42").Value = "'5.51"
$ws.Range("E42").Value = "  -1.46%  "

# Row 43
$ws.Range("D43").Value = "'60.60"
$ws.Range("E43").Value = "  -4.77%  "

# Row 44
$ws.Range("D44").Value = "'0.198"
$ws.Range("E44").Value = "  +1.56%  "

# Row 45
$ws.Range("E45").Value = "  -0.65%  "

# Row 46
$ws.Range("D46").Value = "'0.0994"
$ws.Range("E46").Value = "  +0.33%  "

# Row 47
$ws.Range("D47").Value = "'99.57"
$ws.Range("E47").Value = "  -1.89%  "

# Row 48
$ws.Range("B48").Value = "WOONetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D48").Value = "'0.446"
$ws.Range("E48").Value = "  +22.25%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  -0.62%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  -0.79%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.14"
$ws.Range("E51").Value = "  -1.83%  "
